# added 4wk low sales check
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 154
$ws.Range("B3").Value = 156
$ws.Range("B4").Value = 161
$ws.Range("B5").Value = 165
$ws.Range("B6").Value = 167
$ws.Range("B7").Value = 219
$ws.Range("B8").Value = 221
$ws.Range("B9").Value = 224
$ws.Range("B10").Value = 226
$ws.Range("B11").Value = 228
$ws.Range("B12").Value = 230
$ws.Range("B13").Value = 233
$ws.Range("B14").Value = 235
$ws.Range("B15").Value = 237
